$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1:F3").NumberFormat = "@"

$headers = @("company_code", "code", "name", "type", "address", "phone")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$row2 = @("C001", "B001", "HEAD OFFICE", "branch", "Jl. Soekarno Hatta", "08771939021")
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

$row3 = @("C001", "B002", "PARTNER", "partner", "Jl. Bung Hatta", "08771238769")
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}
